# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit-tracking tables (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the upstream scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1109.8
$ws.Range("I2").Value = 1349.75
$ws.Range("K2").Value = 1349.75
$ws.Range("M2").Value = -1236.75
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 2599.8572
$ws.Range("I4").Value = 1839.8
$ws.Range("K4").Value = 1839.8
$ws.Range("M4").Value = -1725.8
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 188.625
$ws.Range("I9").Value = 159
$ws.Range("K9").Value = 159
$ws.Range("M9").Value = 10
# Row 16 (Leve Item ID 2146)
$ws.Range("H16").Value = 2742.5
$ws.Range("I16").Value = 2489.2856
$ws.Range("J16").Value = 3333.3333
$ws.Range("K16").Value = 2489.2856
$ws.Range("L16").Value = 3333.3333
$ws.Range("M16").Value = -2259.2856
$ws.Range("N16").Value = -3793.3333
# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 22012.25
$ws.Range("I18").Value = 12219.8
$ws.Range("J18").Value = 38333
$ws.Range("K18").Value = 12219.8
$ws.Range("L18").Value = 38333
$ws.Range("M18").Value = -11935.8
$ws.Range("N18").Value = -38901
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 26750
$ws.Range("J32").Value = 26750
$ws.Range("L32").Value = 26750
$ws.Range("N32").Value = -27402
# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 248.36363
$ws.Range("I55").Value = 179.125
$ws.Range("J55").Value = 433
$ws.Range("K55").Value = 179.125
$ws.Range("L55").Value = 433
$ws.Range("M55").Value = 34.875
$ws.Range("N55").Value = -861
# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 36000
$ws.Range("I74").Value = 36000
$ws.Range("K74").Value = 36000
$ws.Range("M74").Value = -35064
# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 36000
$ws.Range("I77").Value = 36000
$ws.Range("K77").Value = 180000
$ws.Range("M77").Value = -175320
# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 6020.2856
$ws.Range("I100").Value = 5228.8
$ws.Range("K100").Value = 5228.8
$ws.Range("M100").Value = -4687.8
# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 3384.2
$ws.Range("I111").Value = 3384.2
$ws.Range("K111").Value = 10152.6
$ws.Range("M111").Value = -7085.599999999999
# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 15829.25
$ws.Range("I135").Value = 2900
$ws.Range("K135").Value = 26100
$ws.Range("M135").Value = -23565

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 50021260
$ws.Range("I32").Value = 50021260
$ws.Range("K32").Value = 50021260
$ws.Range("M32").Value = -50020973
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 15636651
$ws.Range("I74").Value = 25001206
$ws.Range("K74").Value = 25001206
$ws.Range("M74").Value = -25000332
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 15636651
$ws.Range("I77").Value = 25001206
$ws.Range("K77").Value = 125006030
$ws.Range("M77").Value = -125001662
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1544.6666
$ws.Range("I122").Value = 1206
$ws.Range("K122").Value = 3618
$ws.Range("M122").Value = -1168
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 9478.143
$ws.Range("I132").Value = 5753.5
$ws.Range("J132").Value = 14444.333
$ws.Range("K132").Value = 17260.5
$ws.Range("L132").Value = 43332.999
$ws.Range("M132").Value = -14730.5
$ws.Range("N132").Value = -48392.999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -102
$ws.Range("N22").ClearContents()
# Row 96 (Leve Item ID 19525)
$ws.Range("H96").Value = 36176.668
$ws.Range("I96").Value = 11371
$ws.Range("J96").Value = 70904.60000000001
$ws.Range("K96").Value = 11371
$ws.Range("L96").Value = 70904.60000000001
$ws.Range("M96").Value = -8625
$ws.Range("N96").Value = -76396.60000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10 (Leve Item ID 1997)
$ws.Range("H10").Value = 29095.666
$ws.Range("I10").Value = 3639.5
$ws.Range("J10").Value = 80008
$ws.Range("K10").Value = 3639.5
$ws.Range("L10").Value = 80008
$ws.Range("M10").Value = -3500.5
$ws.Range("N10").Value = -80286
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 2245
$ws.Range("I99").Value = 1304.3334
$ws.Range("J99").Value = 2715.3333
$ws.Range("K99").Value = 1304.3334
$ws.Range("L99").Value = 2715.3333
$ws.Range("M99").Value = 193.6666
$ws.Range("N99").Value = -5711.3333
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1462.625
$ws.Range("I105").Value = 1107.8334
$ws.Range("K105").Value = 1107.8334
$ws.Range("M105").Value = 639.1666
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1602.1
$ws.Range("I107").Value = 1156.4286
$ws.Range("J107").Value = 2642
$ws.Range("K107").Value = 1156.4286
$ws.Range("L107").Value = 2642
$ws.Range("M107").Value = 763.5714
$ws.Range("N107").Value = -6482
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 2245
$ws.Range("I126").Value = 1304.3334
$ws.Range("J126").Value = 2715.3333
$ws.Range("K126").Value = 3913.0002
$ws.Range("L126").Value = 8145.999899999999
$ws.Range("M126").Value = -1443.0002
$ws.Range("N126").Value = -13085.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 123 (Leve Item ID 36037)
$ws.Range("H123").Value = 5115
$ws.Range("J123").Value = 5999.5
$ws.Range("L123").Value = 17998.5
$ws.Range("N123").Value = -22898.5
# Row 125 (Leve Item ID 36043)
$ws.Range("H125").Value = 20508.25
$ws.Range("J125").Value = 20508.25
$ws.Range("L125").Value = 61524.75
$ws.Range("N125").Value = -71364.75
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 5653.125
$ws.Range("J131").Value = 4237.375
$ws.Range("L131").Value = 12712.125
$ws.Range("N131").Value = -22792.125
# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 5165
$ws.Range("I137").Value = 5464.2
$ws.Range("J137").Value = 4666.3335
$ws.Range("K137").Value = 16392.6
$ws.Range("L137").Value = 13999.0005
$ws.Range("M137").Value = -11292.6
$ws.Range("N137").Value = -24199.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11 (Leve Item ID 4422)
$ws.Range("H11").Value = 9498352
$ws.Range("I11").Value = 6326237
$ws.Range("K11").Value = 6326237
$ws.Range("M11").Value = -6326098
# Row 12 (Leve Item ID 4093)
$ws.Range("H12").Value = 2500
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -2860
$ws.Range("N12").Value = -1780
# Row 14 (Leve Item ID 4198)
$ws.Range("H14").Value = 218.75
$ws.Range("I14").Value = 200
$ws.Range("J14").Value = 225
$ws.Range("K14").Value = 200
$ws.Range("L14").Value = 225
$ws.Range("M14").Value = -32
$ws.Range("N14").Value = -561
# Row 15 (Leve Item ID 12018)
$ws.Range("H15").Value = 62998.5
$ws.Range("J15").Value = 62998.5
$ws.Range("L15").Value = 62998.5
$ws.Range("N15").Value = -63574.5
# Row 81 (Leve Item ID 12018)
$ws.Range("H81").Value = 62998.5
$ws.Range("J81").Value = 62998.5
$ws.Range("L81").Value = 62998.5
$ws.Range("N81").Value = -64994.5
# Row 84 (Leve Item ID 12018)
$ws.Range("H84").Value = 62998.5
$ws.Range("J84").Value = 62998.5
$ws.Range("L84").Value = 188995.5
$ws.Range("N84").Value = -198979.5
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 4046.6667
$ws.Range("I113").Value = 3938.5715
$ws.Range("J113").Value = 4425
$ws.Range("K113").Value = 3938.5715
$ws.Range("L113").Value = 4425
$ws.Range("M113").Value = -1768.5715
$ws.Range("N113").Value = -8765
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 111114520
$ws.Range("I132").Value = 125003464
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 375010392
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -375007862
$ws.Range("N132").Value = -14060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 30 (Leve Item ID 1688)
$ws.Range("H30").Value = 16
$ws.Range("I30").Value = 16
$ws.Range("K30").Value = 16
$ws.Range("M30").Value = 92
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3444.889
$ws.Range("I40").Value = 2334.4167
$ws.Range("K40").Value = 2334.4167
$ws.Range("M40").Value = -2198.4167
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2866.0386
$ws.Range("J46").Value = 3359
$ws.Range("L46").Value = 3359
$ws.Range("N46").Value = -3735
# Row 80 (Leve Item ID 12027)
$ws.Range("H80").Value = 79564
$ws.Range("J80").Value = 79564
$ws.Range("L80").Value = 79564
$ws.Range("N80").Value = -81810
# Row 83 (Leve Item ID 12027)
$ws.Range("H83").Value = 79564
$ws.Range("J83").Value = 79564
$ws.Range("L83").Value = 238692
$ws.Range("N83").Value = -249924
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 100002080
$ws.Range("I93").Value = 100002080
$ws.Range("K93").Value = 100002080
$ws.Range("M93").Value = -100000832

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 19 (Leve Item ID 2666)
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 3080.65
$ws.Range("I132").Value = 2101.75
$ws.Range("J132").Value = 6996.25
$ws.Range("K132").Value = 6305.25
$ws.Range("L132").Value = 20988.75
$ws.Range("M132").Value = -3775.25
$ws.Range("N132").Value = -26048.75
